$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 760, pushing existing rows 760-864 down to 761-865
$ws.Rows.Item(760).Insert()

# Populate the newly inserted row 760 with the new record
$ws.Cells.Item(760, 1).Value = 5
$ws.Cells.Item(760, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(760, 3).Value = "Maule"
$ws.Cells.Item(760, 4).Value = 45131
$ws.Cells.Item(760, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(760, 5).Value = 7
$ws.Cells.Item(760, 6).Value = 100114001
$ws.Cells.Item(760, 7).Value = "Papa"
$ws.Cells.Item(760, 8).Value = "Asterix"
$ws.Cells.Item(760, 9).Value = "1a (cosecha)"
$ws.Cells.Item(760, 10).Value = 1600
$ws.Cells.Item(760, 11).Value = 16000
$ws.Cells.Item(760, 12).Value = 16000
$ws.Cells.Item(760, 13).Value = 16000
$ws.Cells.Item(760, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(760, 15).Value = "Región del Maule"
$ws.Cells.Item(760, 16).Value = 640
$ws.Cells.Item(760, 17).Value = 25
$ws.Cells.Item(760, 18).Value = "Hortaliza"
